$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Materials sheet: the "Textbook" hyperlink cell text was updated to point at
# the new course repo (display text only - the underlying relationship is
# untouched).
# ---------------------------------------------------------------------------
$wsMaterials = $wb.Worksheets.Item("Materials")
$wsMaterials.Range("B1").Value = "https://github.com/ancestor9/https://github.com/ancestor9/Big-Data-Representtion"

# ---------------------------------------------------------------------------
# Teaching plan sheet: the week-by-week topics were replaced with the new
# Big Data Representation syllabus, and column B was widened to fit the
# longer Korean text.
# ---------------------------------------------------------------------------
$wsPlan = $wb.Worksheets.Item("Teaching plan")
$wsPlan.Range("B2").Value  = " - 데이터 변환과 시각화, 데이터의 종류와 분석방법 이론 예시 "
$wsPlan.Range("B3").Value  = " - Python syntax와 자료 형태"
$wsPlan.Range("B4").Value  = " - Numpy"
$wsPlan.Range("B5").Value  = " - Pandas와 시각화"
$wsPlan.Range("B6").Value  = " - 데이터 종류와 분석방법, 상관관계, 카이제곱검증, ANOVA, 회귀분석, 로지스틱회귀분석"
$wsPlan.Range("B7").Value  = " - 데이터 특성공학(Feature Enginnering) - 실수형과 범주형 변환, survival ship bias"
$wsPlan.Range("B8").Value  = " - 데이터 특성공학(Feature Enginnering) - Target mean, WOE, 이동평균법 등"
$wsPlan.Range("B9").Value  = " - Matplotlib, seaborn - 1"
$wsPlan.Range("B10").Value = " - Matplotlib, seaborn - 2"
$wsPlan.Range("B11").Value = " - 차원축소와 시각화(선형 - PCA)"
$wsPlan.Range("B12").Value = " - 차원축소와 시각화(비선형, t-sne)"
$wsPlan.Range("B13").Value = " - 예측모형과 결과 시각화 - I"
$wsPlan.Range("B14").Value = " - 예측모형과 결과 시각화(Pycaret) - 2"
$wsPlan.Range("B15").Value = " - 딥러닝 representation - I"
$wsPlan.Range("B16").Value = " - 딥러닝 representation - II"

$wsPlan.Columns.Item(2).ColumnWidth = 78.29

# ---------------------------------------------------------------------------
# View state: restore the selections that were active when the file was
# last saved, and make "Teaching plan" the active/visible tab.
# ---------------------------------------------------------------------------
$wsBasic = $wb.Worksheets.Item("Basic information")
$wsBasic.Range("B4").Select()

$wsMaterials.Range("B12").Select()

$wsPlan.Range("B9").Select()
